$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits
#    right after the title (Heading1) paragraph.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2) Insert a new paragraph - "Play Book of Souls for Free -
#    Adventure-themed Video Slot" (bold) - right before the final
#    paragraph (the one that currently holds the AI image prompt).
#    Inserting the raw OOXML at a true mid-paragraph position (not
#    exactly on a paragraph boundary) makes the engine split the
#    paragraph cleanly and keeps the exact run layout we want
#    (a leading empty run followed by the bold run).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$prevPara = $d.Paragraphs($count - 1)
$insertPos = $prevPara.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

$newParaXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Book of Souls for Free - Adventure-themed Video Slot</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertRange.InsertXML($newParaXml)

# ------------------------------------------------------------------
# 3) Replace the text of the last paragraph (the italic AI image
#    prompt) with the new "Join Lara Jones..." copy, keeping its
#    existing (italic) character formatting.
# ------------------------------------------------------------------
$oldText = "Create a dynamic and eye-catching feature image for Book of Souls in cartoon style. The image should feature a happy Maya warrior wearing glasses, as this character is prominent in the game's theme. The warrior should be shown holding the Book of Souls and standing in front of the temple's entrance, surrounded by symbols from the game such as the Mayan masks and totems. Use bright colors and bold lines to make the image pop and convey the excitement and adventure of the game. The image should be designed to capture the attention of potential players and entice them to try out the game."
$newText = "Join Lara Jones on her quest to uncover the treasures of the Book of Souls in this adventure-themed video slot. Play for free and enjoy gameplay features that increase chances of winning big."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
